$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2788.3333
$ws.Range("I32").Value = 1619.25
$ws.Range("J32").Value = 3723.6
$ws.Range("K32").Value = 1619.25
$ws.Range("L32").Value = 3723.6
$ws.Range("M32").Value = -1293.25
$ws.Range("N32").Value = -4375.6

$ws.Range("H34").Value = 6697.4
$ws.Range("I34").Value = 6697.4
$ws.Range("K34").Value = 6697.4
$ws.Range("M34").Value = -6494.4

$ws.Range("H36").Value = 6697.4
$ws.Range("I36").Value = 6697.4
$ws.Range("K36").Value = 6697.4
$ws.Range("M36").Value = -5982.4

$ws.Range("H51").Value = 7073
$ws.Range("I51").Value = 6497.5
$ws.Range("J51").Value = 7456.6665
$ws.Range("K51").Value = 6497.5
$ws.Range("L51").Value = 7456.6665
$ws.Range("M51").Value = -6013.5
$ws.Range("N51").Value = -8424.666499999999

$ws.Range("H88").Value = 1928
$ws.Range("J88").Value = 2240
$ws.Range("L88").Value = 2240
$ws.Range("N88").Value = -3052

$ws.Range("H91").Value = 1928
$ws.Range("J91").Value = 2240
$ws.Range("L91").Value = 2240
$ws.Range("N91").Value = -5048

$ws.Range("H112").Value = 6338452.5
$ws.Range("I112").Value = 1555.5
$ws.Range("J112").Value = 7746651.5
$ws.Range("K112").Value = 4666.5
$ws.Range("L112").Value = 23239954.5
$ws.Range("M112").Value = -3558.5
$ws.Range("N112").Value = -23242170.5

$ws.Range("H135").Value = 2289.5715
$ws.Range("I135").Value = 2286.923
$ws.Range("J135").Value = 2324
$ws.Range("K135").Value = 20582.307
$ws.Range("L135").Value = 20916
$ws.Range("M135").Value = -18047.307
$ws.Range("N135").Value = -25986

$ws.Range("H137").Value = 48439.08
$ws.Range("I137").Value = 60998.95
$ws.Range("K137").Value = 182996.85
$ws.Range("M137").Value = -180446.85

$ws.Range("H138").Value = 3104.33
$ws.Range("I138").Value = 2127.4443
$ws.Range("J138").Value = 3465.6438
$ws.Range("K138").Value = 6382.3329
$ws.Range("L138").Value = 10396.9314
$ws.Range("M138").Value = -1242.3329
$ws.Range("N138").Value = -20676.9314

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10378.6455
$ws.Range("I32").Value = 9015.053
$ws.Range("K32").Value = 9015.053
$ws.Range("M32").Value = -8728.053

$ws.Range("H61").Value = 5636.6763
$ws.Range("J61").Value = 9949.857
$ws.Range("L61").Value = 9949.857
$ws.Range("N61").Value = -10373.857

$ws.Range("H74").Value = 96964.336
$ws.Range("I74").Value = 96964.336
$ws.Range("K74").Value = 96964.336
$ws.Range("M74").Value = -96090.336

$ws.Range("H77").Value = 96964.336
$ws.Range("I77").Value = 96964.336
$ws.Range("K77").Value = 484821.68
$ws.Range("M77").Value = -480453.68

$ws.Range("H97").Value = 2008.3636
$ws.Range("I97").Value = 1516.1765
$ws.Range("J97").Value = 3681.8
$ws.Range("K97").Value = 1516.1765
$ws.Range("L97").Value = 3681.8
$ws.Range("M97").Value = -1020.1765
$ws.Range("N97").Value = -4673.8

$ws.Range("H122").Value = 5846.189
$ws.Range("J122").Value = 13294.75
$ws.Range("L122").Value = 39884.25
$ws.Range("N122").Value = -44784.25

$ws.Range("H132").Value = 5784.1763
$ws.Range("I132").Value = 3312
$ws.Range("K132").Value = 9936
$ws.Range("M132").Value = -7406

$ws.Range("H133").Value = 69261
$ws.Range("J133").Value = 69261
$ws.Range("L133").Value = 69261
$ws.Range("N133").Value = -74321

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 5636.6763
$ws.Range("J136").Value = 9949.857
$ws.Range("L136").Value = 29849.571
$ws.Range("N136").Value = -34949.571

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3625.1765
$ws.Range("I134").Value = 4163.4287
$ws.Range("K134").Value = 12490.2861
$ws.Range("M134").Value = -9955.286100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1686.88
$ws.Range("I16").Value = 1576.3
$ws.Range("K16").Value = 1576.3
$ws.Range("M16").Value = -1289.3

$ws.Range("H81").Value = 44999
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 44999
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H99").Value = 1113483.4
$ws.Range("I99").Value = 1668743.6
$ws.Range("K99").Value = 1668743.6
$ws.Range("M99").Value = -1667245.6

$ws.Range("H113").Value = 1686.88
$ws.Range("I113").Value = 1576.3
$ws.Range("K113").Value = 1576.3
$ws.Range("M113").Value = 593.7

$ws.Range("H126").Value = 1113483.4
$ws.Range("I126").Value = 1668743.6
$ws.Range("K126").Value = 5006230.800000001
$ws.Range("M126").Value = -5003760.800000001

$ws.Range("H132").Value = 3647.1667
$ws.Range("I132").Value = 3415.1365
$ws.Range("K132").Value = 10245.4095
$ws.Range("M132").Value = -7715.4095

$ws.Range("H134").Value = 4283.7295
$ws.Range("I134").Value = 4844.433
$ws.Range("K134").Value = 14533.299
$ws.Range("M134").Value = -11998.299

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1804
$ws.Range("I5").Value = 1804
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 5412
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -5300
$ws.Range("N5").ClearContents()

$ws.Range("H39").Value = 2307.875
$ws.Range("I39").Value = 949
$ws.Range("J39").Value = 2760.8333
$ws.Range("K39").Value = 2847
$ws.Range("L39").Value = 8282.499899999999
$ws.Range("M39").Value = -2553
$ws.Range("N39").Value = -8870.499899999999

$ws.Range("H51").Value = 17299
$ws.Range("I51").Value = 17299
$ws.Range("K51").Value = 51897
$ws.Range("M51").Value = -51437

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws.Range("H105").Value = 19000
$ws.Range("J105").Value = 19000
$ws.Range("L105").Value = 57000
$ws.Range("N105").Value = -62242

$ws.Range("H121").Value = 2864
$ws.Range("J121").Value = 3385.5557
$ws.Range("L121").Value = 10156.6671
$ws.Range("N121").Value = -12776.6671

$ws.Range("H129").Value = 2406.8462
$ws.Range("I129").Value = 2294.2
$ws.Range("J129").Value = 2477.25
$ws.Range("K129").Value = 6882.599999999999
$ws.Range("L129").Value = 7431.75
$ws.Range("M129").Value = -1882.599999999999
$ws.Range("N129").Value = -17431.75

$ws.Range("H133").Value = 3500
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

$ws.Range("H134").Value = 9080.5
$ws.Range("I134").Value = 9080.5
$ws.Range("K134").Value = 27241.5
$ws.Range("M134").Value = -22171.5

$ws.Range("H135").Value = 1804
$ws.Range("I135").Value = 1804
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 16236
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -13701
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 3035.5
$ws.Range("I136").Value = 1311.5454
$ws.Range("K136").Value = 3934.6362
$ws.Range("M136").Value = 1165.3638

$ws.Range("H137").Value = 2647.5
$ws.Range("I137").Value = 1863.3334
$ws.Range("J137").Value = 5000
$ws.Range("K137").Value = 5590.0002
$ws.Range("L137").Value = 15000
$ws.Range("M137").Value = -490.0002000000004
$ws.Range("N137").Value = -25200

$ws.Range("H138").Value = 7776.6665
$ws.Range("I138").Value = 7776.6665
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 23329.9995
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -18189.9995
$ws.Range("N138").ClearContents()

$ws.Range("H139").Value = 2954.5
$ws.Range("I139").Value = 1959.4
$ws.Range("J139").Value = 3507.3333
$ws.Range("K139").Value = 5878.200000000001
$ws.Range("L139").Value = 10521.9999
$ws.Range("M139").Value = -738.2000000000007
$ws.Range("N139").Value = -20801.9999

$ws.Range("H140").Value = 4631266
$ws.Range("I140").Value = 17857824
$ws.Range("K140").Value = 53573472
$ws.Range("M140").Value = -53568292

$ws.Range("H141").Value = 1793.1111
$ws.Range("I141").Value = 1579.75
$ws.Range("K141").Value = 4739.25
$ws.Range("M141").Value = 440.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 134471.67
$ws.Range("I132").Value = 172657.83
$ws.Range("K132").Value = 517973.49
$ws.Range("M132").Value = -515443.49

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3014.121
$ws.Range("J136").Value = 2997.4167
$ws.Range("L136").Value = 8992.250100000001
$ws.Range("N136").Value = -14092.2501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2670.423
$ws.Range("I132").Value = 2538.7273
$ws.Range("K132").Value = 7616.1819
$ws.Range("M132").Value = -5086.1819
